## Splits the literal "<exp>...</exp>" abbreviation-expansion markers that
## currently live inline inside plain-text runs into their own
## Courier-New / gray / small-caps "tag" runs, matching the look used
## elsewhere in the transcription for markers like <lb/>, <del>, <page>, etc.
##
## Each of the six target runs has the shape:
##     PREFIX + "<exp>" + MID + "</exp>" + SUFFIX
## and must become up to five runs:
##     PREFIX              (unchanged formatting)
##     "<exp>"             (Courier New / a9a9a9 / 7pt)
##     MID                 (unchanged formatting)
##     "</exp>"            (Courier New / a9a9a9 / 7pt)
##     SUFFIX               (unchanged formatting, omitted if empty)

$d = $word.ActiveDocument

function Format-TagRange([int]$rs, [int]$re, [string]$txt) {
    # Re-stamp the [rs, re) range (which already holds $txt) with the
    # gray Courier-New "markup tag" look used throughout the document.
    $sub = $d.Range($rs, $re)
    $sub.Select()
    $sel = $word.Selection
    $sel.Find.ClearFormatting()
    $sel.Find.Replacement.ClearFormatting()
    $sel.Find.Replacement.Font.Name = "Courier New"
    $sel.Find.Replacement.Font.NameFarEast = "Courier New"
    $sel.Find.Replacement.Font.NameBi = "Courier New"
    $sel.Find.Replacement.Font.Color = 11119017
    $sel.Find.Replacement.Font.Size = 7
    $sel.Find.Execute($txt, $true, $false, $false, $false, $false, $true, 1, $false, $txt, 2)
}

function Split-ExpRun([string]$prefix, [string]$mid, [string]$suffix) {
    $full = $prefix + "<exp>" + $mid + "</exp>" + $suffix

    $finder = $d.Content
    $finder.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $finder.Find.Found) {
        Write-Output "NOT FOUND: $full"
        return
    }

    $p1 = $finder.Start
    $p2 = $p1 + $prefix.Length
    $p3 = $p2 + 5                      # "<exp>".Length
    $p4 = $p3 + $mid.Length
    $p5 = $p4 + 6                      # "</exp>".Length
    $p6 = $p5 + $suffix.Length

    Format-TagRange $p2 $p3 "<exp>"
    Format-TagRange $p4 $p5 "</exp>"
}

Split-ExpRun " de mesme ligue que le canon co" "mm" "e"
Split-ExpRun "que chasque refoulem" "ent" ""
Split-ExpRun " done adva" "n" "tage"
Split-ExpRun "Elle acompaigne fort de poincte en blanc la grand colevrine pourceq" "ue" ""
Split-ExpRun "pans de longueur co" "mm" "e est la vache de "
Split-ExpRun "pieces on donne renfort a la culasse co" "mm" "e de trois balles Elles peuvent"

Write-Output "all done"
